$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.809.20"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.317.63"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'97.21"
$ws.Range("E5").Value = "  +5.83%  "
$ws.Range("D6").Value = "'272.66"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "'0.627"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "'45.32"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("D11").Value = "'0.0951"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  -4.32%  "
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "2.658.44"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").Value = "'15.52"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "'0.876"
$ws.Range("E16").Value = "  +8.03%  "
$ws.Range("D17").Value = "2.320.60"
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").Value = "43.765.83"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +3.81%  "
$ws.Range("D20").Value = "'6.41"
$ws.Range("E20").Value = "  +4.71%  "
$ws.Range("D21").Value = "'73.38"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("D22").Value = "'239.86"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").Value = "'2.27"
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("D24").Value = "'9.41"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Value = "'11.37"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("D30").Value = "'38.21"
$ws.Range("E30").Value = "  -7.19%  "
$ws.Range("D31").Value = "'22.42"
$ws.Range("E31").Value = "  +6.66%  "
$ws.Range("D32").Value = "'175.16"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").Value = "'0.0919"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("E37").Value = "  -4.04%  "
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").Value = "'3.39"
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("D40").Value = "'0.246"
$ws.Range("E40").Value = "  +8.28%  "
$ws.Range("E41").Value = "  +9.61%  "
$ws.Range("D42").Value = "'1.41"
$ws.Range("E42").Value = "  +21.72%  "
$ws.Range("D43").Value = "'12.32"
$ws.Range("E43").Value = "  -5.60%  "
$ws.Range("D44").Value = "'62.75"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("E45").Value = "  +9.96%  "
$ws.Range("D46").Value = "'5.36"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("E47").Value = "  +3.49%  "
$ws.Range("D48").Value = "'100.48"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "'0.194"
$ws.Range("E50").Value = "  +17.12%  "
$ws.Range("D51").Value = "2.545.53"
$ws.Range("E51").Value = "  +2.68%  "

# Strip the quote-prefix style iron_native applies for forced-text numeric-looking
# values above, so the cells keep style index 0 (matching the original, unstyled cells)
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
